$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly price record as row 333, shifting the existing
# rows 333:387 down to 334:388 (dimension grows from A1:R387 to A1:R388).
$ws.Rows(333).Insert()

$ws.Range("A333").Value = 7
$ws.Range("B333").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C333").Value = "Ñuble"
$ws.Range("D333").Value = 45218
$ws.Range("E333").Value = 16
$ws.Range("F333").Value = 100112043
$ws.Range("G333").Value = "Pepino ensalada"
$ws.Range("H333").Value = "Sin especificar"
$ws.Range("I333").Value = "Primera"
$ws.Range("J333").Value = 100
$ws.Range("K333").Value = 14000
$ws.Range("L333").Value = 14000
$ws.Range("M333").Value = 14000
$ws.Range("N333").Value = "$/caja 60 unidades"
$ws.Range("O333").Value = "Región de Arica y Parinacota"
$ws.Range("P333").Value = 233
$ws.Range("Q333").Value = 60
$ws.Range("R333").Value = "Hortaliza"
